$wb = $excel.ActiveWorkbook

# --- Struct sheet ("Struct" = sheet6.xml): add new rows 18-20 ---
$wsStruct = $wb.Worksheets.Item("Struct")

$wsStruct.Range("A18").Value = "Orbiter body h"
$wsStruct.Range("B18").Value = 1

$wsStruct.Range("A19").Value = "Orbiter body w"
$wsStruct.Range("B19").Value = 1

$wsStruct.Range("A20").Value = "Orbiter body l"
$wsStruct.Range("B20").Value = 1

# update the saved selection on the Struct sheet (it stays not-active)
$wsStruct.Range("J23").Select()

# --- Prop sheet ("Prop" = sheet5.xml): update B24, add rows 32-33 ---
$wsProp = $wb.Worksheets.Item("Prop")

$wsProp.Range("B24").Value = 200

$wsProp.Range("A32").Value = "Maintenance thruster mass"
$wsProp.Range("B32").Value = 0.40033000000000002
$wsProp.Range("C32").Value = "kg"

$wsProp.Range("A33").Value = "Circularisation propulsion system dry mass"
$wsProp.Range("B33").Value = 82.84502106177365
$wsProp.Range("D33").Value = "kg"

# Prop becomes the active sheet/tab, with A33:D33 selected
$wsProp.Activate()
$wsProp.Range("A33:D33").Select()
